$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source feed was re-synced one day later: every existing row's
# "Förändrad" (changed) date in column C moves from 45181 to 45182
# (rows 2 through 428 inclusive).
for ($r = 2; $r -le 428; $r++) {
    $ws.Cells.Item($r, 3).Value = 45182
}

# Row 428 picks up an explicit row height now that a new row follows it.
$ws.Rows.Item(428).RowHeight = 15

# A brand new logging notification was appended as row 429.
$ws.Cells.Item(429, 1).Value = "A 42726-2023"
$ws.Cells.Item(429, 2).Value = 45181
$ws.Cells.Item(429, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(429, 3).Value = 45182
$ws.Cells.Item(429, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(429, 4).Value = "SÖDERMANLANDS LÄN"
$ws.Cells.Item(429, 5).Value = "STRÄNGNÄS"
$ws.Cells.Item(429, 7).Value = 2.7
$ws.Cells.Item(429, 8).Value = 0
$ws.Cells.Item(429, 9).Value = 0
$ws.Cells.Item(429, 10).Value = 0
$ws.Cells.Item(429, 11).Value = 0
$ws.Cells.Item(429, 12).Value = 0
$ws.Cells.Item(429, 13).Value = 0
$ws.Cells.Item(429, 14).Value = 0
$ws.Cells.Item(429, 15).Value = 0
$ws.Cells.Item(429, 16).Value = 0
$ws.Cells.Item(429, 17).Value = 0

# Column R holds a wrap-text note column; keep it consistent with the
# rest of the table (empty, wrap-text style) for the new row.
$ws.Cells.Item(429, 18).WrapText = $true
